# Selected another equivalent parts to cover shortage
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: J7 (SMD,P=1.27mm, JTAG) - swap JLCPCB part# C920875 -> C2935953
$ws.Range("D23").Value = "C2935953"

# Row 28: C21 (1206, 4.7nF 1kV) - swap JLCPCB part# C377102 -> C106074
$ws.Range("D28").Value = "C106074"

# Update the last active/selected cell to D29
$ws.Range("D29").Select()
